# WorkListSample.xlsx - edc-patcher "robustness" tweak:
# insert a new "S4 System Name" column right after "AI2 Site Reference"
# (i.e. before the old "NGR" column), fill in its two values, size the
# column to fit its content, and leave the selection where the author
# left it (cell E7, no frozen/scrolled top-left cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column C - everything from the old C onward (NGR,
# Install Date, AI2 Equipment SAI/PLI Number, ...) shifts right by one.
$ws.Columns("C:C").Insert()

# Populate the header + sample row for the newly inserted column.
$ws.Range("C1").Value = "S4 System Name"
$ws.Range("C2").Value = "EA Event Duration Monitoring"

# Size the new column to fit its (now longest) content, same as the
# other "bestFit" columns on this sheet.
$ws.Columns("C:C").ColumnWidth = 27.7109375

# Match the author's final selection/view state.
[void]$ws.Range("E7").Select()
